$d = $word.ActiveDocument

# --- Paragraph 1 formatting: add a paragraph border (space-only, no line)
# and change the left indent from 120 -> 225 twips (11.25 pt). ---
$p1 = $d.Paragraphs(1)

$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5

$p1.Format.LeftIndent = 11.25

# --- Text change: update the bookmark id text and drop the trailing
# space run that used to follow it. ---
$d.Content.Find.Execute("**ID__AFFARS_pgi_5317_topic_13__ID**", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_SMC_PGI_5317__ID**", 2) | Out-Null

$p1 = $d.Paragraphs(1)
$pEnd = $p1.Range.End
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}
